$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 224.82353
$ws.Range("J33").Value = 499
$ws.Range("L33").Value = 499
$ws.Range("N33").Value = -957
$ws.Range("H53").Value = 20834596
$ws.Range("I53").Value = 55556476
$ws.Range("J53").Value = 1469.7
$ws.Range("K53").Value = 55556476
$ws.Range("L53").Value = 1469.7
$ws.Range("M53").Value = -55555839
$ws.Range("N53").Value = -2743.7
$ws.Range("H69").Value = 17499.75
$ws.Range("J69").Value = 20000
$ws.Range("L69").Value = 60000
$ws.Range("N69").Value = -61748
$ws.Range("H72").Value = 17499.75
$ws.Range("J72").Value = 20000
$ws.Range("L72").Value = 180000
$ws.Range("N72").Value = -188736
$ws.Range("H98").Value = 919.069
$ws.Range("I98").Value = 616.8148
$ws.Range("K98").Value = 616.8148
$ws.Range("M98").Value = 881.1852
$ws.Range("H122").Value = 919.069
$ws.Range("I122").Value = 616.8148
$ws.Range("K122").Value = 1850.4444
$ws.Range("M122").Value = 599.5556000000001
$ws.Range("H132").Value = 2498.4
$ws.Range("I132").Value = 2404.6
$ws.Range("K132").Value = 7213.799999999999
$ws.Range("M132").Value = -4683.799999999999
$ws.Range("H137").Value = 2161.2424
$ws.Range("I137").Value = 2089.08
$ws.Range("K137").Value = 6267.24
$ws.Range("M137").Value = -3717.24

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2959.19
$ws.Range("I32").Value = 2959.19
$ws.Range("K32").Value = 2959.19
$ws.Range("M32").Value = -2672.19
$ws.Range("H44").Value = 52032.668
$ws.Range("J44").Value = 52032.668
$ws.Range("L44").Value = 52032.668
$ws.Range("N44").Value = -53008.668
$ws.Range("H61").Value = 3794.9
$ws.Range("I61").Value = 4118.625
$ws.Range("K61").Value = 4118.625
$ws.Range("M61").Value = -3906.625
$ws.Range("H74").Value = 1417.4
$ws.Range("I74").Value = 1469.6666
$ws.Range("K74").Value = 1469.6666
$ws.Range("M74").Value = -595.6666
$ws.Range("H77").Value = 1417.4
$ws.Range("I77").Value = 1469.6666
$ws.Range("K77").Value = 7348.333000000001
$ws.Range("M77").Value = -2980.333000000001
$ws.Range("H125").Value = 88988.55
$ws.Range("J125").Value = 88988.55
$ws.Range("L125").Value = 88988.55
$ws.Range("N125").Value = -98828.55
$ws.Range("H132").Value = 4440.2173
$ws.Range("I132").Value = 4414.773
$ws.Range("K132").Value = 13244.319
$ws.Range("M132").Value = -10714.319
$ws.Range("H136").Value = 3794.9
$ws.Range("I136").Value = 4118.625
$ws.Range("K136").Value = 12355.875
$ws.Range("M136").Value = -9805.875

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 47329.39
$ws.Range("I134").Value = 4026.182
$ws.Range("K134").Value = 12078.546
$ws.Range("M134").Value = -9543.545999999998

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 336.2143
$ws.Range("I22").Value = 361.36365
$ws.Range("J22").Value = 244
$ws.Range("K22").Value = 361.36365
$ws.Range("L22").Value = 244
$ws.Range("M22").Value = -11.36365000000001
$ws.Range("N22").Value = -944
$ws.Range("H31").Value = 21740.195
$ws.Range("I31").Value = 1367.6666
$ws.Range("J31").Value = 44659.293
$ws.Range("K31").Value = 1367.6666
$ws.Range("L31").Value = 44659.293
$ws.Range("M31").Value = -1072.6666
$ws.Range("N31").Value = -45249.293
$ws.Range("H34").Value = 21740.195
$ws.Range("I34").Value = 1367.6666
$ws.Range("J34").Value = 44659.293
$ws.Range("K34").Value = 1367.6666
$ws.Range("L34").Value = 44659.293
$ws.Range("M34").Value = -1165.6666
$ws.Range("N34").Value = -45063.293
$ws.Range("H99").Value = 5041.1665
$ws.Range("I99").Value = 4394.9
$ws.Range("J99").Value = 5849
$ws.Range("K99").Value = 4394.9
$ws.Range("L99").Value = 5849
$ws.Range("M99").Value = -2896.9
$ws.Range("N99").Value = -8845
$ws.Range("H126").Value = 5041.1665
$ws.Range("I126").Value = 4394.9
$ws.Range("J126").Value = 5849
$ws.Range("K126").Value = 13184.7
$ws.Range("L126").Value = 17547
$ws.Range("M126").Value = -10714.7
$ws.Range("N126").Value = -22487
$ws.Range("H132").Value = 2544.8333
$ws.Range("I132").Value = 2163.7273
$ws.Range("J132").Value = 3143.7144
$ws.Range("K132").Value = 6491.1819
$ws.Range("L132").Value = 9431.143199999999
$ws.Range("M132").Value = -3961.1819
$ws.Range("N132").Value = -14491.1432
$ws.Range("H137").Value = 45000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 837294.75
$ws.Range("I80").Value = 671110
$ws.Range("K80").Value = 671110
$ws.Range("M80").Value = -670112
$ws.Range("H83").Value = 837294.75
$ws.Range("I83").Value = 671110
$ws.Range("K83").Value = 3355550
$ws.Range("M83").Value = -3350558
$ws.Range("H132").Value = 54958.05
$ws.Range("I132").Value = 5459.222
$ws.Range("J132").Value = 500447.5
$ws.Range("K132").Value = 16377.666
$ws.Range("L132").Value = 1501342.5
$ws.Range("M132").Value = -13847.666
$ws.Range("N132").Value = -1506402.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 12598.5
$ws.Range("J42").Value = 12598.5
$ws.Range("L42").Value = 12598.5
$ws.Range("N42").Value = -13724.5
$ws.Range("H46").Value = 1630.6316
$ws.Range("I46").Value = 1676
$ws.Range("K46").Value = 1676
$ws.Range("M46").Value = -1488
$ws.Range("H49").Value = 12598.5
$ws.Range("J49").Value = 12598.5
$ws.Range("L49").Value = 12598.5
$ws.Range("N49").Value = -12892.5
$ws.Range("H122").Value = 3584.6667
$ws.Range("I122").Value = 3666.6924
$ws.Range("K122").Value = 11000.0772
$ws.Range("M122").Value = -8550.0772
$ws.Range("H136").Value = 232227.3
$ws.Range("I136").Value = 326354.3
$ws.Range("K136").Value = 979062.8999999999
$ws.Range("M136").Value = -976512.8999999999
$ws.Range("H139").Value = 56166.668
$ws.Range("J139").Value = 56166.668
$ws.Range("L139").Value = 56166.668
$ws.Range("N139").Value = -66446.66800000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 27000
$ws.Range("I42").Value = 27000
$ws.Range("K42").Value = 27000
$ws.Range("M42").Value = -26622
$ws.Range("H75").Value = 32997.5
$ws.Range("H78").Value = 32997.5
$ws.Range("H122").Value = 27780480
$ws.Range("I122").Value = 40002250
$ws.Range("K122").Value = 120006750
$ws.Range("M122").Value = -120004300
$ws.Range("H126").Value = 1870.8
$ws.Range("I126").Value = 1812.5
$ws.Range("J126").Value = 2104
$ws.Range("K126").Value = 5437.5
$ws.Range("L126").Value = 6312
$ws.Range("M126").Value = -2967.5
$ws.Range("N126").Value = -11252
$ws.Range("H132").Value = 20037.75
$ws.Range("I132").Value = 2573.2827
$ws.Range("J132").Value = 77421
$ws.Range("K132").Value = 7719.848100000001
$ws.Range("L132").Value = 232263
$ws.Range("M132").Value = -5189.848100000001
$ws.Range("N132").Value = -237323
$ws.Range("H136").Value = 56783.406
$ws.Range("I136").Value = 2446.1562
$ws.Range("J136").Value = 404541.8
$ws.Range("K136").Value = 7338.4686
$ws.Range("L136").Value = 1213625.4
$ws.Range("M136").Value = -4788.4686
$ws.Range("N136").Value = -1218725.4
$ws.Range("H139").Value = 99990
$ws.Range("J139").Value = 99990
$ws.Range("L139").Value = 99990
$ws.Range("N139").Value = -110270
